$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts old "Tipo" column to E)
$ws.Columns("D").Insert()

# New header for the inserted column
# (Columns.Insert() already carries over the bold/centered header style
#  from the old column D / new column E, matching the rest of row 1)
$ws.Range("D1").Value = "MAE"

# New MAE values
$ws.Range("D2").Value = 0.6472536726993832
$ws.Range("D3").Value = 0.3447233267553575
